$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.430.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.952.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.18%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.840"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.238.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("E16").Value = "  +1.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.950.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.382.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0852"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("E22").Value = "  +2.49%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "

$ws.Range("E25").Value = "  +3.46%  "

$ws.Range("E26").Value = "  +7.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.77%  "

$ws.Range("E31").Value = "  +1.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0609"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("E34").Value = "  +6.86%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  +2.65%  "

$ws.Range("E37").Value = "  +4.56%  "

$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.45%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.360.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.133.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
